$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Mise a jour de la data table "stock" (quantites disponibles) pour le nouveau workflow
$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 28
$ws.Range("B6").Value = 128
